# Add a new "Spain" market test-data sheet, cloned from the existing
# "Italy" sheet (same layout/styles), with Spain-specific values.

$wb = $excel.ActiveWorkbook
$italy = $wb.Worksheets.Item("Italy")

# Duplicate "Italy" to create "Spain" right after it, then rename.
$italy.Copy($null, $italy) | Out-Null
$spain = $wb.Worksheets.Item($wb.Worksheets.Count)
$spain.Name = "Spain"

# Fill in the Spain-specific market name and product/test code.
$spain.Range("B2").Value = "Spain Market"
$spain.Range("B4").Value = "NGC-3103/T2045"

# Re-fit columns/rows for the new text (mirrors Excel auto-fit after edit).
$spain.Columns.Item(1).ColumnWidth = 24.276042
$spain.Columns.Item(2).ColumnWidth = 14.385417
$spain.Columns.Item(3).ColumnWidth = 10.498698
$spain.Columns.Item(4).ColumnWidth = 20.94401
$spain.Rows.Item(3).RowHeight = 28.8
$spain.Rows.Item(4).RowHeight = 28.8

# Italy is no longer the active tab - reset its selection to the full table.
$italy.Range("A1:D12").Select() | Out-Null

# Spain becomes the active tab, with D9 selected.
$spain.Range("D9").Select() | Out-Null
